$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Renumber existing "Work" items 57 -> 58, 56 -> 57, ... 39 -> 40
#    (processed from highest to lowest so we never double-increment)
#    This makes room for the brand-new item 39 that we insert below.
# ------------------------------------------------------------------
for ($n = 57; $n -ge 39; $n--) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Style = "WorkNumber"
    $target = [string]$n
    $replacement = [string]($n + 1)
    $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2) | Out-Null
}

# ------------------------------------------------------------------
# 2. Locate the last German entry ("... Babel, 2024.") which is
#    immediately followed by an empty "Work" spacer paragraph. The
#    new item 39 (plus its own spacer) needs to be inserted right
#    after that existing spacer paragraph, before the blank
#    separator / "Translated from the Portuguese" heading.
# ------------------------------------------------------------------
$anchorRng = $d.Content
$anchorRng.Find.ClearFormatting()
$found = $anchorRng.Find.Execute("Babel, 2024.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find insertion anchor paragraph"
}
$anchorPara = $anchorRng.Paragraphs(1)
$spacerPara = $anchorPara.Next()

$insertionPoint = $spacerPara.Range.Duplicate
$insertionPoint.Collapse(0)

$snippet = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Work"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"></w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Work"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r>
  <w:r><w:rPr><w:rStyle w:val="WorkNumber"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">39</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r>
  <w:r><w:rPr><w:rStyle w:val="BookTitle"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Zehn Millionen Kinder [School for Barbarians]</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">by</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Erika Mann.</w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Afarsemon, 2024.</w:t></w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Work"/>
    <w:rPr><w:lang w:val="en-US"/></w:rPr>
  </w:pPr>
  <w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"></w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint.InsertXML($snippet)

Write-Host "Done. Total paragraphs now:" $d.Paragraphs.Count
